$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.159.23"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.70%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.612.21"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.27%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.07%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'559.60"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +5.05%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'143.97"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.72%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.23%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.598"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +5.35%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'6.81"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -1.26%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.30%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +5.49%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.11%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.074.67"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.38%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'59.101.94"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.75%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.26%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.621.05"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.09%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.08%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'4.46"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.37%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'337.78"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.78%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.05%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.30%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.09%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'66.08"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.36%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  +3.84%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.59%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.994"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.59%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'7.20"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.36%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'0.0₃0767"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +4.34%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.13%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.36%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'6.04"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +3.89%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'154.83"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +2.46%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.30%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +1.75%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.913"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +11.11%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.907"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +9.24%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.11%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.42%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +3.09%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'3.62"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.33%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'285.22"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +1.27%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.23%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.602"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.49%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.0542"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.69%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0959"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +2.33%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'10.61"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -1.03%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +4.19%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +2.08%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.955.74"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.78%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'117.97"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +5.68%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'18.12"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.83%  "
$ws.Range('E51').Style = 'Normal'
